# Update cosinor analysis results (row 2 and row 3) with re-run values
# from the CircadiPy sawtooth_10 simulation analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("G2").Value = [double]"0.0004874112272066355"
$ws.Range("H2").Value = [double]"0.001485253369957398"
$ws.Range("K2").Value = [double]"5.188429808816737"
$ws.Range("L2").Value = "[2.2443792668563614, 8.132480350777113]"
$ws.Range("M2").Value = [double]"0.0006141406263899807"
$ws.Range("N2").Value = [double]"0.0006141406263899807"
$ws.Range("P2").Value = "[-2.251631971942233, -0.7421580242714629]"
$ws.Range("Q2").Value = [double]"0.0001218549013350945"
$ws.Range("R2").Value = [double]"0.0001218549013350945"
$ws.Range("S2").Value = [double]"10.85250687009681"
$ws.Range("T2").Value = "[9.04596208190457, 12.659051658289055]"
$ws.Range("X2").Value = [double]"3.069889889889965"
$ws.Range("Y2").Value = [double]"9.313733733733955"

# --- Row 3 ---
$ws.Range("E3").Value = [double]"24.75000000000043"
$ws.Range("G3").Value = [double]"2.933180706843341e-05"
$ws.Range("H3").Value = [double]"0.0003134235367074911"
$ws.Range("I3").Value = [double]"1.110223024625157e-16"
$ws.Range("K3").Value = [double]"5.176938540239227"
$ws.Range("L3").Value = "[2.5847318552512313, 7.7691452252272235]"
$ws.Range("M3").Value = [double]"0.0001062708519370403"
$ws.Range("N3").Value = [double]"0.0002125417038740807"
$ws.Range("O3").Value = [double]"-2.993789996213697"
$ws.Range("P3").Value = "[-3.6101585248459287, -2.377421467581465]"
$ws.Range("S3").Value = [double]"9.914401291579935"
$ws.Range("T3").Value = "[8.347947776447558, 11.480854806712312]"
$ws.Range("W3").Value = [double]"11.792792792793"
$ws.Range("X3").Value = [double]"9.364864864865027"
$ws.Range("Y3").Value = [double]"14.22072072072097"
